# Adjust load up settings
# - Updates the controls_array text used by rows 18-33 and 43-51 (shared string)
# - Updates heating_min_shed (P) and percent_shed_met (Q) for rows 27-33 and 43-49

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Update the controls_array string (column D) for the affected rows.
# Old: "N, N, N, N, N, N, N, N, N, N, N, N, L, L, L, L, S, S, S, S, N, N, N, N"
# New: "N, N, N, N, N, N, N, N, N, N, L, L, L, L, L, L, S, S, S, S, N, N, N, N"
$newControlsArray = "N, N, N, N, N, N, N, N, N, N, L, L, L, L, L, L, S, S, S, S, N, N, N, N"

$controlRows = 18..33 + 43..51
foreach ($r in $controlRows) {
    $ws.Range("D$r").Value = $newControlsArray
}

# Update heating_min_shed (P) / percent_shed_met (Q) pairs for the rows that changed.
$updates = @{
    27 = @{ P = 15;  Q = 0.9375 }
    28 = @{ P = 12;  Q = 0.95 }
    29 = @{ P = 55;  Q = 0.7708333333333334 }
    30 = @{ P = 69;  Q = 0.7125 }
    31 = @{ P = 94;  Q = 0.6083333333333334 }
    32 = @{ P = 101; Q = 0.5791666666666666 }
    33 = @{ P = 65;  Q = 0.7291666666666667 }
    43 = @{ P = 13;  Q = 0.9458333333333333 }
    44 = @{ P = 121; Q = 0.4958333333333333 }
    45 = @{ P = 99;  Q = 0.5875 }
    46 = @{ P = 50;  Q = 0.7916666666666666 }
    47 = @{ P = 80;  Q = 0.6666666666666667 }
    48 = @{ P = 129; Q = 0.4625 }
    49 = @{ P = 164; Q = 0.3166666666666667 }
}

foreach ($r in $updates.Keys) {
    $ws.Range("P$r").Value = $updates[$r].P
    $ws.Range("Q$r").Value = $updates[$r].Q
}
